$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows for the "origin" / "dest" metadata entries, right
# after the arr_time row (row 3) and before the carrier row (old row 4).
$ws.Rows("4:5").Insert()

# New row 4: origin -> origem
$ws.Range("A4").Value = "nyflights"
$ws.Range("B4").Value = "origin"
$ws.Range("C4").Value = "string"
$ws.Range("D4").Value = "origem"
$ws.Range("E4").Value = "string"
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0.05
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0

# New row 5: dest -> destino
$ws.Range("A5").Value = "nyflights"
$ws.Range("B5").Value = "dest"
$ws.Range("C5").Value = "string"
$ws.Range("D5").Value = "destino"
$ws.Range("E5").Value = "string"
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0.05
$ws.Range("H5").Value = 1
$ws.Range("I5").Value = 0

# Relabel the "str" type marker as "string" across the remaining rows
# (dep_time, arr_time, carrier, flight, tailnum) except tipo_original for
# tailnum, which keeps its original "str" label.
$ws.Range("C2").Value = "string"
$ws.Range("C3").Value = "string"
$ws.Range("C6").Value = "string"
$ws.Range("E6").Value = "string"
$ws.Range("C7").Value = "string"
$ws.Range("E7").Value = "string"
$ws.Range("E8").Value = "string"

# The distance row's I cell (previously the last row, carrying the
# underline placeholder style) loses that styling once it is no longer
# the table's trailing edge.
$ws.Range("I10").Font.Underline = -4142

# Update selection/active cell to match the saved view state.
$ws.Range("E14").Select()
